$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shape = $s.Shapes.Item($i)
    if ($shape.Name -eq "Рисунок 5") {
        $shape.Delete()
    }
}
